# The diff adds two new, empty, centered+bold paragraphs right after the
# title paragraph ("FUNCIONALIDADES PRINCIPALES DE CROPSCAN") and before
# the first bullet ("El aplicativo dara...").  Each new paragraph looks
# like:
#   <w:p>
#     <w:pPr>
#       <w:jc w:val="center"/>
#       <w:rPr><w:b/><w:bCs/></w:rPr>
#     </w:pPr>
#   </w:p>
# i.e. no run at all - just paragraph mark formatting (centered, bold).
#
# A plain InsertParagraphAfter() on the title paragraph's range correctly
# picks up the jc/b/bCs paragraph-mark formatting, but it also manifests
# as a stray empty <w:r><w:rPr>...</w:rPr></w:r> run in the new paragraph,
# which the diff does not contain. Instead we build the exact OOXML for
# the two new paragraphs and drop it in with Range.InsertXML, which
# inserts the given well-formed paragraphs verbatim (no synthesized run).

$d = $word.ActiveDocument

$titlePara = $d.Paragraphs(1)
$insertionPoint = $d.Range($titlePara.Range.End, $titlePara.Range.End)

$emptyCenteredBoldParagraph = '<w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr></w:p>'

$package = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
            '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
                '<w:body>' + $emptyCenteredBoldParagraph + $emptyCenteredBoldParagraph + '</w:body>' +
            '</w:document>' +
        '</pkg:xmlData>' +
    '</pkg:part>' +
'</pkg:package>'

[void]$insertionPoint.InsertXML($package)
